# This script re-applies the latest cryptocurrency price/volume snapshot
# (values scraped from coinranking.com) onto columns D (Price) and E
# (Volume(1h)) of the active worksheet, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "23.444.76"; E = "+1.84%" },
    @{ Row = 3; D = "1.629.10"; E = "+2.67%" },
    @{ Row = 4; D = "0.9967"; E = "-0.60%" },
    @{ Row = 5; D = "307.50"; E = "+2.10%" },
    @{ Row = 6; D = "0.9971"; E = "-0.54%" },
    @{ Row = 7; D = "0.3780"; E = "+0.37%" },
    @{ Row = 8; D = "53.11"; E = "+5.26%" },
    @{ Row = 9; D = "0.3665"; E = "+2.15%" },
    @{ Row = 10; D = "1.286"; E = "+5.52%" },
    @{ Row = 11; D = "0.08196"; E = "+2.00%" },
    @{ Row = 12; D = "0.9970"; E = "-0.60%" },
    @{ Row = 13; D = $null; E = "+6.11%" },
    @{ Row = 14; D = "6.677"; E = "+3.04%" },
    @{ Row = 15; D = "7.467"; E = "+1.73%" },
    @{ Row = 16; D = "0.00001265"; E = "+3.48%" },
    @{ Row = 17; D = "1.630.91"; E = "+2.67%" },
    @{ Row = 18; D = $null; E = "+2.76%" },
    @{ Row = 19; D = "0.06931"; E = "+2.13%" },
    @{ Row = 20; D = $null; E = "+2.93%" },
    @{ Row = 21; D = "6.592"; E = "+2.49%" },
    @{ Row = 22; D = "0.9985"; E = "-0.42%" },
    @{ Row = 23; D = "13.00"; E = "+1.52%" },
    @{ Row = 24; D = "23.462.45"; E = "+1.96%" },
    @{ Row = 25; D = "3.125"; E = "+13.02%" },
    @{ Row = 26; D = "2.437"; E = "+3.01%" },
    @{ Row = 27; D = $null; E = "+2.89%" },
    @{ Row = 28; D = "150.92"; E = "+2.40%" },
    @{ Row = 29; D = "5.307"; E = "+1.93%" },
    @{ Row = 30; D = "136.43"; E = "+2.58%" },
    @{ Row = 31; D = "2.432"; E = "+4.88%" },
    @{ Row = 32; D = "6.949"; E = "+6.49%" },
    @{ Row = 33; D = "1.808.20"; E = "+2.53%" },
    @{ Row = 34; D = "0.9790"; E = "+4.02%" },
    @{ Row = 35; D = "0.02808"; E = "+4.96%" },
    @{ Row = 36; D = "10.49"; E = "+4.57%" },
    @{ Row = 37; D = "0.07501"; E = "+2.17%" },
    @{ Row = 38; D = "6.256"; E = "+2.99%" },
    @{ Row = 39; D = $null; E = "+2.32%" },
    @{ Row = 40; D = "0.08847"; E = "+0.98%" },
    @{ Row = 41; D = "1.411"; E = "+5.48%" },
    @{ Row = 42; D = "0.7174"; E = "+4.21%" },
    @{ Row = 43; D = $null; E = "+7.62%" },
    @{ Row = 44; D = "16.18"; E = "+8.89%" },
    @{ Row = 45; D = "0.6635"; E = "+3.61%" },
    @{ Row = 46; D = $null; E = "+5.17%" },
    @{ Row = 47; D = $null; E = "+1.19%" },
    @{ Row = 48; D = "0.9959"; E = "-0.54%" },
    @{ Row = 49; D = "0.08023"; E = "+1.78%" },
    @{ Row = 50; D = "132.22"; E = "+0.62%" },
    @{ Row = 51; D = $null; E = "+1.50%" }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        # Force text formatting so price strings such as "0.9967" or
        # "23.444.76" are preserved verbatim and not auto-converted to numbers.
        $ws.Range("D$($u.Row)").NumberFormat = "@"
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = "  $($u.E)  "
}
